$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ C = 16.911; D = 5486.6 }
    3  = @{ C = 17.295; D = 191.4 }
    4  = @{ C = 18.119; D = 175.8 }
    5  = @{ C = 17.086; D = 176.4 }
    6  = @{ C = 17.516; D = 182.8 }
    7  = @{ C = 16.667; D = 184 }
    8  = @{ C = 15.976; D = 123 }
    9  = @{ C = 15.696; D = 177.6 }
    10 = @{ C = 15.684; D = 154.2 }
    11 = @{ C = 17.083; D = 154.2 }
    12 = @{ C = 16.188; D = 154.4 }
    13 = @{ C = 17.337; D = 176.4 }
    14 = @{ C = 17.594; D = 183.8 }
    15 = @{ C = 17.273; D = 123 }
    16 = @{ C = 17.613; D = 133.2 }
}

foreach ($row in $data.Keys) {
    $ws.Cells.Item($row, 3).Value = $data[$row].C
    $ws.Cells.Item($row, 4).Value = $data[$row].D
}
